$d = $word.ActiveDocument

# 1. "HashSet<User> onlineUsers" -> "HashMap<String, User> userNameToUser"
$d.Content.Find.Execute("HashSet<User> onlineUsers", $true, $false, $false, $false, $false,
                         $true, 1, $false, "HashMap<String, User> userNameToUser", 2)

# 2. "conversationNameToConversations" -> "conversationNameToConversation"
$d.Content.Find.Execute("conversationNameToConversations", $true, $false, $false, $false, $false,
                         $true, 1, $false, "conversationNameToConversation", 2)
